# Team member report push
# Update the per-member "Role duties and work to be performed next week",
# "Issues encountered" / "Issues resolved" cells, and append extra work
# performed this week for Brody and Vasilis.

$wb = $excel.ActiveWorkbook

# --- Arpit (sheet "Arpit") ---
$ws = $wb.Worksheets.Item("Arpit")
$ws.Range("B8").Value = "Development Team"
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 36
$ws.Activate()
$ws.Range("B9").Select()

# --- Brody ---
$ws = $wb.Worksheets.Item("Brody")
$ws.Range("B7").Value = "Updated product backlog; Updated SRS, URN documents; Participated in Usability Study, created forms for it; Create software architecture and design pattern document; Updated SRS document; Added Use case 34; PDFized several documents"
$ws.Range("B8").Value = "Development Team"
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 144
$ws.Activate()
$ws.Range("B8").Select()

# --- Michael ---
$ws = $wb.Worksheets.Item("Michael")
$ws.Range("B8").Value = "Scrum Master"
$ws.Range("B9").Value = "Out of Town 4/11-4/16"
$ws.Range("A7").WrapText = $true
$ws.Activate()
$ws.Range("B14").Select()

# --- Sakshyam ---
$ws = $wb.Worksheets.Item("Sakshyam")
$ws.Range("B8").Value = "Development Team"
$ws.Range("B9").Value = "Unfamiliarity with Unit Testing"
$ws.Range("B10").Value = "N/A"
$ws.Range("A7").WrapText = $true
$ws.Activate()
$ws.Range("B10").Select()

# --- Vasilis ---
$ws = $wb.Worksheets.Item("Vasilis")
$ws.Range("B7").Value = "Release new version APK for deliverable 3; Updated URN and SRS documents; Participated in Usability Study, created script, recorded and took note of activities; Updated Kanban board.; Created misuse case diagram for use case 34; Reworded privacy notification; Completed acceptance & unit testing; Combined elements into test plan document"
$ws.Range("B8").Value = "Development Team"
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 198
$ws.Activate()
$ws.Range("A7").Select()

# --- Yong ---
$ws = $wb.Worksheets.Item("Yong")
$ws.Range("B8").Value = "Product Owner"
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 144
$ws.Activate()
$ws.Range("A14").Select()

# Leave Vasilis as the active/visible sheet, matching the saved view state.
$vasilis = $wb.Worksheets.Item("Vasilis")
$vasilis.Activate()
$vasilis.Range("A7").Select()
